# "Main.xlsx" / sheet "Rules": cell B11 is updated from the text "R40" to
# the text "1". This introduces a new entry ("1") in the shared-strings
# table; the cell's existing formatting (style) is left alone.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A leading apostrophe tells Excel to store the value as literal text
# rather than re-interpreting the digit "1" as the number 1, matching the
# target workbook where B11 keeps its string ("s") cell type.
$ws.Range("B11").Value = "'1"
